$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.252.12'
$ws.Range('E2').Value = '  +1.48%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.908.38'
$ws.Range('E3').Value = '  +2.12%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.83'
$ws.Range('E5').Value = '  +0.93%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.11%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5251'
$ws.Range('E7').Value = '  +3.37%  '

# Row 8
$ws.Range('E8').Value = '  +3.63%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07301'
$ws.Range('E9').Value = '  +1.62%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.31'
$ws.Range('E10').Value = '  +3.04%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9007'
$ws.Range('E11').Value = '  +0.93%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.933.05'
$ws.Range('E12').Value = '  +3.45%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07674'
$ws.Range('E13').Value = '  +2.00%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.09'
$ws.Range('E14').Value = '  +0.37%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.261'
$ws.Range('E15').Value = '  +0.79%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.08%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008665'
$ws.Range('E17').Value = '  +2.04%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.50'
$ws.Range('E18').Value = '  +2.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9991'
$ws.Range('E19').Value = '  -0.17%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.307.41'
$ws.Range('E20').Value = '  +1.52%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.086'
$ws.Range('E21').Value = '  +1.53%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.147.42'
$ws.Range('E22').Value = '  +1.30%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.65'
$ws.Range('E23').Value = '  +2.84%  '

# Row 24
$ws.Range('E24').Value = '  +1.03%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.313'
$ws.Range('E25').Value = '  +10.49%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.93'
$ws.Range('E26').Value = '  -1.43%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.16'
$ws.Range('E27').Value = '  +1.64%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.736'
$ws.Range('E28').Value = '  -2.64%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.83'
$ws.Range('E29').Value = '  +1.27%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.957'
$ws.Range('E30').Value = '  +4.78%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.818'
$ws.Range('E31').Value = '  +2.57%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09223'
$ws.Range('E32').Value = '  +0.88%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05085'
$ws.Range('E33').Value = '  +0.36%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7975'
$ws.Range('E34').Value = '  +6.74%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.249'
$ws.Range('E35').Value = '  +8.17%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.987'
$ws.Range('E36').Value = '  +0.21%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.303'
$ws.Range('E37').Value = '  +2.29%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.625'
$ws.Range('E38').Value = '  +3.83%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5683'
$ws.Range('E39').Value = '  +1.64%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01996'
$ws.Range('E40').Value = '  +0.11%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.075'
$ws.Range('E41').Value = '  +0.08%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.033'
$ws.Range('E42').Value = '  +5.14%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.649'
$ws.Range('E43').Value = '  +0.64%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.95'
$ws.Range('E44').Value = '  +3.11%  '

# Row 45
$ws.Range('E45').Value = '  +3.25%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4868'
$ws.Range('E46').Value = '  +2.85%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.25'
$ws.Range('E47').Value = '  +1.02%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9992'
$ws.Range('E48').Value = '  -0.11%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.605'
$ws.Range('E49').Value = '  +2.59%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.47'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.23'
$ws.Range('E51').Value = '  +1.93%  '
